# Updated names in Jan sheet
#
# - Switch the active sheet from "March" to "January".
# - Enter the four team member names into January!B3:B6 (Excel will
#   automatically back them with shared strings).
# - Autofit column B so it displays the full names, and leave the
#   selection on B7 (the cell right below the last entry), matching
#   what Excel does after typing a list of values and pressing Enter.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("January")
$ws.Activate()

$ws.Range("B3").Value = "Nitesh"
$ws.Range("B4").Value = "Gautami"
$ws.Range("B5").Value = "Pratiksha"
$ws.Range("B6").Value = "Pruthviraj"

$ws.Columns("B").AutoFit() | Out-Null

$ws.Range("B7").Select()
